$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.975.10'
$ws.Range("E2").Value = '  -0.45%  '

$ws.Range("D3").Value = '1.866.94'
$ws.Range("E3").Value = '  -2.89%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '''317.90'
$ws.Range("E5").Value = '  -2.90%  '

$ws.Range("E6").Value = '  +0.03%  '

$ws.Range("D7").Value = '''0.5077'
$ws.Range("E7").Value = '  -1.70%  '

$ws.Range("D8").Value = '''0.3915'
$ws.Range("E8").Value = '  -2.31%  '

$ws.Range("D9").Value = '''0.08130'
$ws.Range("E9").Value = '  -3.91%  '

$ws.Range("D10").Value = '''41.96'
$ws.Range("E10").Value = '  -2.17%  '

$ws.Range("D11").Value = '''1.087'
$ws.Range("E11").Value = '  -3.01%  '

$ws.Range("D12").Value = '''22.65'
$ws.Range("E12").Value = '  +6.47%  '

$ws.Range("D13").Value = '1.857.55'
$ws.Range("E13").Value = '  -3.12%  '

$ws.Range("D14").Value = '''6.242'
$ws.Range("E14").Value = '  -1.44%  '

$ws.Range("D15").Value = '''7.149'
$ws.Range("E15").Value = '  -2.58%  '

$ws.Range("E16").Value = '  +0.12%  '

$ws.Range("D17").Value = '''91.38'
$ws.Range("E17").Value = '  -3.56%  '

$ws.Range("D18").Value = '''0.00001073'
$ws.Range("E18").Value = '  -3.87%  '

$ws.Range("D19").Value = '''0.06355'
$ws.Range("E19").Value = '  -5.64%  '

$ws.Range("E20").Value = '  -1.36%  '

$ws.Range("E21").Value = '  +0.09%  '

$ws.Range("D22").Value = '29.966.76'
$ws.Range("E22").Value = '  -0.51%  '

$ws.Range("D23").Value = '''5.774'
$ws.Range("E23").Value = '  -4.82%  '

$ws.Range("D24").Value = '''11.04'
$ws.Range("E24").Value = '  -1.48%  '

$ws.Range("E25").Value = '  -0.06%  '

$ws.Range("D26").Value = '2.084.83'
$ws.Range("E26").Value = '  -2.55%  '

$ws.Range("D27").Value = '''160.55'
$ws.Range("E27").Value = '  -0.24%  '

$ws.Range("D28").Value = '''20.79'
$ws.Range("E28").Value = '  -0.93%  '

$ws.Range("D29").Value = '''2.212'
$ws.Range("E29").Value = '  -10.14%  '

$ws.Range("D30").Value = '''126.19'
$ws.Range("E30").Value = '  -2.01%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = '''1.042'
$ws.Range("E31").Value = '  -3.26%  '

$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").Value = '''0.1029'
$ws.Range("E32").Value = '  -2.99%  '

$ws.Range("D33").Value = '''5.852'
$ws.Range("E33").Value = '  -3.58%  '

$ws.Range("D34").Value = '''3.728'
$ws.Range("E34").Value = '  +1.62%  '

$ws.Range("D35").Value = '''0.02407'
$ws.Range("E35").Value = '  -4.18%  '

$ws.Range("B36").Value = 'InternetComputer(DFINITY)'
$ws.Range("C36").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D36").Value = '''5.176'
$ws.Range("E36").Value = '  -0.56%  '

$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").Value = '''0.06314'
$ws.Range("E37").Value = '  -4.22%  '

$ws.Range("D38").Value = '''0.2130'
$ws.Range("E38").Value = '  -4.30%  '

$ws.Range("D39").Value = '''1.169'
$ws.Range("E39").Value = '  -5.57%  '

$ws.Range("D40").Value = '''8.463'
$ws.Range("E40").Value = '  -6.01%  '

$ws.Range("D41").Value = '''0.6254'
$ws.Range("E41").Value = '  -4.41%  '

$ws.Range("E42").Value = '  -2.86%  '

$ws.Range("D43").Value = '''11.24'
$ws.Range("E43").Value = '  -1.29%  '

$ws.Range("D44").Value = '''0.9996'
$ws.Range("E44").Value = '  +0.03%  '

$ws.Range("D45").Value = '''0.5854'
$ws.Range("E45").Value = '  -4.61%  '

$ws.Range("D46").Value = '''12.78'
$ws.Range("E46").Value = '  -3.32%  '

$ws.Range("D47").Value = '''3.617'
$ws.Range("E47").Value = '  -3.88%  '

$ws.Range("D48").Value = '''1.980'
$ws.Range("E48").Value = '  -3.60%  '

$ws.Range("D49").Value = '''121.95'
$ws.Range("E49").Value = '  -3.13%  '

$ws.Range("D50").Value = '''1.200'
$ws.Range("E50").Value = '  -3.42%  '

$ws.Range("D51").Value = '''1.151'
$ws.Range("E51").Value = '  -0.34%  '
